# Job para atualizar pautas
#
# 1) Shrink the two "pauta" table's left indent / cell margins by a hair
#    and nudge a few column widths (1408->1407, 4394->4393, 1248->1250).
# 2) Make the header-row label text (AUTORIA, O QUE FAZ, POSICAO SAL,
#    PRIORITARIOS, COMISSAO / PLENARIO, PAUTA) white like the PROJETO
#    header cell already is.
# 3) Merge the trailing "Senado Federal"+" " runs into a single run.

$d = $word.ActiveDocument

# Merge the "Senado Federal" + " " runs into one run, keeping bold.
# (Done before touching $d.Tables -- walking the Tables/Cells collection
# first throws off Paragraphs.Item() indexing in this host.)
for ($pi = 1; $pi -le $d.Paragraphs.Count; $pi++) {
    $p = $d.Paragraphs.Item($pi)
    if ($p.Range.Text -eq "Senado Federal " + [char]13) {
        $start = $p.Range.Start
        $end = $p.Range.End - 1
        $r = $d.Range($start, $end)
        $r.Text = ""
        $ins = $d.Range($start, $start)
        $ins.InsertAfter("Senado Federal ")
        $d.Range($start, $start + 15).Font.Bold = $true
    }
}

for ($ti = 1; $ti -le $d.Tables.Count; $ti++) {
    $t = $d.Tables.Item($ti)

    # table-wide indent: -10dxa (-0.5pt) -> -15dxa (-0.75pt)
    $t.Rows.LeftIndent = -0.75

    # table-wide default cell left margin: 98dxa (4.9pt) -> 93dxa (4.65pt)
    $t.LeftPadding = 4.65

    for ($ri = 1; $ri -le $t.Rows.Count; $ri++) {
        $row = $t.Rows.Item($ri)

        # column widths (dxa/20 = pt): col1 1408->1407, col3 4394->4393, col7 1248->1250
        $row.Cells.Item(1).Width = 70.35
        $row.Cells.Item(3).Width = 219.65
        $row.Cells.Item(7).Width = 62.5

        for ($ci = 1; $ci -le $row.Cells.Count; $ci++) {
            $cell = $row.Cells.Item($ci)
            # every cell's explicit left margin override: 98dxa -> 93dxa
            $cell.LeftPadding = 4.65
        }
    }

    # header row label cells (everything right of "PROJETO") go white
    $headerRow = $t.Rows.Item(1)
    for ($ci = 2; $ci -le $headerRow.Cells.Count; $ci++) {
        $headerRow.Cells.Item($ci).Range.Font.Color = 16777215
    }
}
